# Corrección graficas y preguntas: actualizar los datos de conteo (Si/No)
# y dejar la selección activa en A3, tal como quedó el libro tras la edición.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Corregir los valores de la tabla de conteo
$ws.Range("A2").Value = 53
$ws.Range("B2").Value = 0

# Dejar seleccionada la celda A3 (como quedó guardado en el libro original)
$ws.Activate()
$ws.Range("A3").Select()

# Reflejar la nueva posición de la ventana del libro (best effort; algunas
# propiedades de geometría de ventana pueden no persistir en el runtime)
try {
    $win = $excel.Windows.Item(1)
    $win.Left = 1500
    $win.Top = 1500
} catch {
}
